$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '24.175.05'
Set-TextCell 2 5 '  -2.50%  '
Set-TextCell 3 4 '1.643.72'
Set-TextCell 3 5 '  -2.48%  '
Set-TextCell 4 4 '1.000'
Set-TextCell 4 5 '  -0.51%  '
Set-TextCell 5 4 '307.84'
Set-TextCell 5 5 '  -2.03%  '
Set-TextCell 6 4 '1.0000'
Set-TextCell 6 5 '  -0.38%  '
Set-TextCell 7 4 '0.3892'
Set-TextCell 7 5 '  -0.99%  '
Set-TextCell 8 4 '0.3863'
Set-TextCell 8 5 '  -2.67%  '
Set-TextCell 9 5 '  -0.53%  '
Set-TextCell 10 4 '49.65'
Set-TextCell 10 5 '  -4.53%  '
Set-TextCell 11 5 '  -4.93%  '
Set-TextCell 12 4 '0.08645'
Set-TextCell 12 5 '  -0.45%  '
Set-TextCell 13 4 '23.64'
Set-TextCell 13 5 '  -6.33%  '
Set-TextCell 14 4 '7.119'
Set-TextCell 14 5 '  -2.77%  '
Set-TextCell 15 4 '0.00001290'
Set-TextCell 15 5 '  -2.37%  '
Set-TextCell 16 4 '7.457'
Set-TextCell 16 5 '  -4.42%  '
Set-TextCell 17 4 '1.646.06'
Set-TextCell 17 5 '  +0.62%  '
Set-TextCell 18 4 '94.86'
Set-TextCell 18 5 '  +0.56%  '
Set-TextCell 19 4 '0.06902'
Set-TextCell 19 5 '  -2.93%  '
Set-TextCell 20 4 '20.42'
Set-TextCell 20 5 '  +1.33%  '
Set-TextCell 21 4 '6.901'
Set-TextCell 21 5 '  -3.50%  '
Set-TextCell 22 4 '1.000'
Set-TextCell 22 5 '  -0.39%  '
Set-TextCell 23 4 '13.57'
Set-TextCell 23 5 '  -3.83%  '
Set-TextCell 24 4 '24.170.19'
Set-TextCell 24 5 '  -2.52%  '
Set-TextCell 25 4 '2.387'
Set-TextCell 25 5 '  -0.22%  '
Set-TextCell 26 4 '2.799'
Set-TextCell 26 5 '  +0.77%  '
Set-TextCell 27 4 '22.40'
Set-TextCell 27 5 '  -5.64%  '
Set-TextCell 28 4 '157.86'
Set-TextCell 28 5 '  -2.72%  '
Set-TextCell 29 4 '8.598'
Set-TextCell 29 5 '  +9.60%  '
Set-TextCell 30 4 '140.37'
Set-TextCell 30 5 '  -6.28%  '
Set-TextCell 31 4 '5.362'
Set-TextCell 31 5 '  -7.01%  '
Set-TextCell 32 4 '2.410'
Set-TextCell 32 5 '  -7.56%  '
Set-TextCell 33 4 '1.830.29'
Set-TextCell 33 5 '  +2.20%  '
Set-TextCell 34 4 '7.011'
Set-TextCell 34 5 '  +0.79%  '
Set-TextCell 35 4 '0.08074'
Set-TextCell 35 5 '  -4.68%  '
Set-TextCell 36 4 '0.02908'
Set-TextCell 36 5 '  -5.76%  '
Set-TextCell 37 4 '0.2681'
Set-TextCell 37 5 '  -4.61%  '
Set-TextCell 38 4 '0.9506'
Set-TextCell 38 5 '  -6.10%  '
Set-TextCell 39 4 '0.09203'
Set-TextCell 39 5 '  -3.64%  '
Set-TextCell 40 2 'TrustWalletToken'
Set-TextCell 40 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 40 4 '1.464'
Set-TextCell 40 5 '  +0.41%  '
Set-TextCell 41 2 'FraxShare'
Set-TextCell 41 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 41 4 '9.988'
Set-TextCell 41 5 '  -4.86%  '
Set-TextCell 42 4 '0.7547'
Set-TextCell 42 5 '  -5.21%  '
Set-TextCell 43 4 '13.02'
Set-TextCell 43 5 '  -5.21%  '
Set-TextCell 44 4 '16.05'
Set-TextCell 44 5 '  -3.99%  '
Set-TextCell 45 4 '0.6913'
Set-TextCell 45 5 '  -3.76%  '
Set-TextCell 46 4 '2.460'
Set-TextCell 46 5 '  -4.89%  '
Set-TextCell 47 4 '4.090'
Set-TextCell 47 5 '  -2.45%  '
Set-TextCell 48 4 '0.9997'
Set-TextCell 48 5 '  -0.34%  '
Set-TextCell 49 4 '0.08408'
Set-TextCell 49 5 '  -3.82%  '
Set-TextCell 50 5 '  -5.30%  '
Set-TextCell 51 4 '133.34'
Set-TextCell 51 5 '  -3.51%  '
